$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 12.15509999999999
$ws.Range("E3").Value = 13.4468
$ws.Range("E5").Value = 13.04929999999999
$ws.Range("C9").Value = -11.89530000000001
$ws.Range("E11").Value = 13.48529999999999
$ws.Range("E12").Value = 13.05229999999999
$ws.Range("C13").Value = -11.94419999999999
$ws.Range("C16").Value = -11.8464
$ws.Range("C18").Value = -14.33299999999999
$ws.Range("C20").Value = -13.48559999999999
$ws.Range("E21").Value = 12.89449999999999

$wb.Save()
